$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F15").Value = 50
$ws1.Range("F18").Value = 6576
$ws1.Range("F20").Value = 7380
$ws1.Range("F22").Value = 56464
$ws1.Range("F23").Value = 4536
$ws1.Range("F25").Value = 887
$ws1.Range("F28").Value = 890
$ws1.Range("F35").Value = 1247
$ws1.Range("F36").Value = 1254
$ws1.Range("F39").Value = 199

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F9").Value = 9348
$ws3.Range("F16").Value = 367

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F18").Value = 50
$ws4.Range("F21").Value = 56464
$ws4.Range("F23").Value = 4536
$ws4.Range("F35").Value = 1247
$ws4.Range("F37").Value = 367
